$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain Text (the source data stores price strings
# like "41.384.99" / "6.60" / "5.90" as literal text, not numbers) so the
# trailing zeros and thousand-dot grouping survive the round-trip.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.384.99'
$ws.Range("E2").Value = '  +3.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.255.17'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.68'
$ws.Range("E5").Value = '  +3.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.65'
$ws.Range("E6").Value = '  +4.96%  '
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("E9").Value = '  +3.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.87'
$ws.Range("E10").Value = '  +7.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.05'
$ws.Range("E11").Value = '  +7.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0793'
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.60'
$ws.Range("E14").Value = '  +2.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.603.48'
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.14'
$ws.Range("E16").Value = '  +3.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.239.48'
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.304.11'
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.16'
$ws.Range("E20").Value = '  +8.52%  '
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.90'
$ws.Range("E22").Value = '  +2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.86'
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.42'
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("E25").Value = '  +4.82%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.86'
$ws.Range("E27").Value = '  +3.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.74'
$ws.Range("E28").Value = '  +5.66%  '
$ws.Range("E29").Value = '  +6.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.66'
$ws.Range("E30").Value = '  +5.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.71'
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.70'
$ws.Range("E32").Value = '  +8.01%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.17'
$ws.Range("E34").Value = '  +6.12%  '
$ws.Range("E35").Value = '  +4.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.04'
$ws.Range("E36").Value = '  +8.36%  '
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.61'
$ws.Range("E38").Value = '  +9.29%  '
$ws.Range("E39").Value = '  +2.81%  '
$ws.Range("E40").Value = '  +5.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("E42").Value = '  +7.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.22'
$ws.Range("E43").Value = '  +16.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.064.00'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0276'
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.14'
$ws.Range("E46").Value = '  +5.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.97'
$ws.Range("E47").Value = '  +12.49%  '
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.474.01'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("E50").Value = '  +3.16%  '
$ws.Range("E51").Value = '  +3.50%  '
